$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B39").Value = "Character Info"
$ws.Range("B39:G39").Merge()

$ws.Range("B40").Value = "Name"
$ws.Range("B40:C40").Merge()
$ws.Range("D40").Value = "Color"
$ws.Range("D40:E40").Merge()
$ws.Range("F40").Value = "Row"
$ws.Range("G40").Value = "Col"
